$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 14:50"

$ws.Range("B4").Value = 36249
$ws.Range("C4").Value = 15362
$ws.Range("D4").Value = 16164
$ws.Range("E4").Value = 4723

$ws.Range("B6").Value = 5625
$ws.Range("C6").Value = 531
$ws.Range("D6").Value = 4926
$ws.Range("E6").Value = 168

$ws.Range("B7").Value = 4224
$ws.Range("C7").Value = 3098
$ws.Range("D7").Value = 4612

$ws.Range("B8").Value = 3555
$ws.Range("C8").Value = 425
$ws.Range("D8").Value = 2860
$ws.Range("E8").Value = 270

$ws.Range("B10").Value = 2972
$ws.Range("D10").Value = 2507
$ws.Range("E10").Value = 171

$ws.Range("A12").Value = "Alacant/Alicante"
$ws.Range("B12").Value = 2559
$ws.Range("C12").Value = 214
$ws.Range("D12").Value = 2100
$ws.Range("E12").Value = 245

$ws.Range("A13").Value = "Araba/Alava"
$ws.Range("B13").Value = 2539
$ws.Range("C13").Value = 3098
$ws.Range("D13").Value = 4612
$ws.Range("E13").Value = 176

$ws.Range("A14").Value = "La Rioja"
$ws.Range("B14").Value = 2405
$ws.Range("C14").Value = 843
$ws.Range("D14").Value = 1434
$ws.Range("E14").Value = 128

$ws.Range("A15").Value = "Albacete"
$ws.Range("B15").Value = 2386
$ws.Range("C15").Value = 492
$ws.Range("D15").Value = 7028
$ws.Range("E15").Value = 183

$ws.Range("B16").Value = 2342
$ws.Range("C16").Value = 331
$ws.Range("D16").Value = 1819
$ws.Range("E16").Value = 192

$ws.Range("C18").Value = 162
$ws.Range("D18").Value = 1633
$ws.Range("E18").Value = 110

$ws.Range("C20").Value = 53
$ws.Range("D20").Value = 1490
$ws.Range("E20").Value = 75

$ws.Range("A21").Value = "Salamanca"
$ws.Range("B21").Value = 1579
$ws.Range("C21").Value = 332
$ws.Range("D21").Value = 1069
$ws.Range("E21").Value = 178

$ws.Range("A22").Value = "Pontevedra"
$ws.Range("B22").Value = 1536
$ws.Range("C22").Value = 333
$ws.Range("D22").Value = 1411
$ws.Range("E22").Value = 30

$ws.Range("A23").Value = "Asturias"
$ws.Range("B23").Value = 1522
$ws.Range("C23").Value = 154
$ws.Range("D23").Value = 1292
$ws.Range("E23").Value = 76

$ws.Range("B24").Value = 1424
$ws.Range("C24").Value = 3098
$ws.Range("D24").Value = 4612

$ws.Range("A25").Value = "Cantabria"
$ws.Range("B25").Value = 1384
$ws.Range("C25").Value = 99
$ws.Range("D25").Value = 1217
$ws.Range("E25").Value = 68

$ws.Range("A26").Value = "Valladolid"
$ws.Range("B26").Value = 1352
$ws.Range("C26").Value = 439
$ws.Range("D26").Value = 788
$ws.Range("E26").Value = 125

$ws.Range("A27").Value = "Granada"
$ws.Range("B27").Value = 1340
$ws.Range("C27").Value = 47
$ws.Range("D27").Value = 1188
$ws.Range("E27").Value = 105

$ws.Range("B28").Value = 1317
$ws.Range("C28").Value = 85
$ws.Range("D28").Value = 1053
$ws.Range("E28").Value = 179

$ws.Range("A29").Value = "Leon"
$ws.Range("B29").Value = 1204
$ws.Range("C29").Value = 398
$ws.Range("D29").Value = 647
$ws.Range("E29").Value = 159

$ws.Range("A30").Value = "Murcia"
$ws.Range("B30").Value = 1188
$ws.Range("C30").Value = 113
$ws.Range("D30").Value = 1024
$ws.Range("E30").Value = 51

$ws.Range("B31").Value = 1051
$ws.Range("C31").Value = 311
$ws.Range("D31").Value = 638
$ws.Range("E31").Value = 102

$ws.Range("B32").Value = 964
$ws.Range("C32").Value = 323
$ws.Range("D32").Value = 545
$ws.Range("E32").Value = 96

$ws.Range("A33").Value = "Tenerife"
$ws.Range("B33").Value = 946
$ws.Range("C33").Value = 123
$ws.Range("D33").Value = 1564
$ws.Range("E33").Value = 51

$ws.Range("A34").Value = "Aragon"
$ws.Range("B34").Value = 907
$ws.Range("C34").Value = 29
$ws.Range("D34").Value = 838
$ws.Range("E34").Value = 40

$ws.Range("C35").Value = 55
$ws.Range("D35").Value = 792
$ws.Range("E35").Value = 30

$ws.Range("C36").Value = 39
$ws.Range("D36").Value = 782
$ws.Range("E36").Value = 46

$ws.Range("C38").Value = 62
$ws.Range("D38").Value = 717
$ws.Range("E38").Value = 27

$ws.Range("A39").Value = "Castello/Castellon"
$ws.Range("B39").Value = 787
$ws.Range("C39").Value = 56
$ws.Range("D39").Value = 675
$ws.Range("E39").Value = 56

$ws.Range("A40").Value = "Soria"
$ws.Range("B40").Value = 769
$ws.Range("C40").Value = 136
$ws.Range("D40").Value = 581
$ws.Range("E40").Value = 52

$ws.Range("A41").Value = "Ourense"
$ws.Range("B41").Value = 751
$ws.Range("C41").Value = 333
$ws.Range("D41").Value = 660
$ws.Range("E41").Value = 22

$ws.Range("B42").Value = 662
$ws.Range("C42").Value = 120
$ws.Range("D42").Value = 513
$ws.Range("E42").Value = 29

$ws.Range("B43").Value = 638
$ws.Range("C43").Value = 205
$ws.Range("D43").Value = 364
$ws.Range("E43").Value = 69

$ws.Range("B45").Value = 452
$ws.Range("C45").Value = 104
$ws.Range("D45").Value = 317
$ws.Range("E45").Value = 31

$ws.Range("B46").Value = 435
$ws.Range("E46").Value = 21

$ws.Range("B47").Value = 393
$ws.Range("C47").Value = 44
$ws.Range("D47").Value = 319
$ws.Range("E47").Value = 30

$ws.Range("A49").Value = "Teruel"
$ws.Range("B49").Value = 343
$ws.Range("C49").Value = 65
$ws.Range("D49").Value = 251
$ws.Range("E49").Value = 27

$ws.Range("A50").Value = "Almeria"
$ws.Range("B50").Value = 329
$ws.Range("C50").Value = 37
$ws.Range("D50").Value = 270
$ws.Range("E50").Value = 22

$ws.Range("B51").Value = 323
$ws.Range("C51").Value = 83
$ws.Range("D51").Value = 205
$ws.Range("E51").Value = 35

$ws.Range("C52").Value = 8
$ws.Range("D52").Value = 247
$ws.Range("E52").Value = 11

$ws.Range("A55").Value = "Ceuta"
$ws.Range("B55").Value = 73
$ws.Range("C55").Value = 2
$ws.Range("D55").Value = 70
$ws.Range("E55").Value = 1

$ws.Range("A56").Value = "La Palma"
$ws.Range("B56").Value = 69

$ws.Range("A57").Value = "Lanzarote"
$ws.Range("B57").Value = 68
$ws.Range("C57").Value = 123
$ws.Range("D57").Value = 1564
$ws.Range("E57").Value = 3

$ws.Range("B59").Value = 35
